{"js": "// Insert three new list entries right after the paragraph that ends with\n// \"... du probl\u00e8me des tables dans la bdd\" (25/07 entry: date header,\n// \"Faits :\" sub-header, and the new fact about the double table / matches()\n// issue), matching the existing \"Paragraphedeliste\" numbering (numId=1) at\n// levels 0, 1 and 2 respectively.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph containing the known anchor text.\nconst anchorText = \"R\u00e9solution du probl\u00e8me de connexion de new compte (double hashage) + peut \u00eatre compr\u00e9hension (chatgpt) du probl\u00e8me des tables dans la bdd\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\n// 1) \"25/07 :\" \u2014 top level (ilvl 0) list item, same style/numbering as its\n//    siblings (e.g. \"24/07 :\").\nconst dateHeader = anchor.insertParagraph(\"25/07 :\", Word.InsertLocation.after);\ndateHeader.listItemOrNullObject.level = 0;\nawait context.sync();\n\n// 2) \"Faits :\" \u2014 second level (ilvl 1) list item.\nconst faitsHeader = dateHeader.insertParagraph(\"Faits :\", Word.InsertLocation.after);\nfaitsHeader.listItemOrNullObject.level = 1;\nawait context.sync();\n\n// 3) The new fact itself \u2014 third level (ilvl 2) list item.\nconst factText = \"R\u00e9solution de la double table dans la bdd mais re-probl\u00e8me avec les mdp (.matches())\";\nconst factItem = faitsHeader.insertParagraph(factText, Word.InsertLocation.after);\nfactItem.listItemOrNullObject.level = 2;\nawait context.sync();\n", "ps1": "# Insert three new list entries right after the paragraph that ends with\n# \"... du probl\u00e8me des tables dans la bdd\" (25/07 entry: date header,\n# \"Faits :\" sub-header, and the new fact about the double table / matches()\n# issue), matching the existing \"Paragraphedeliste\" numbering (numId=1) at\n# levels 0, 1 and 2 respectively (Word COM's ListLevelNumber is 1-based,\n# i.e. OOXML ilvl = ListLevelNumber - 1).\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its exact text. Range.Text includes a\n# trailing paragraph-mark character (CR, 13) which must be trimmed before\n# comparing.\n$anchorText = \"R\u00e9solution du probl\u00e8me de connexion de new compte (double hashage) + peut \u00eatre compr\u00e9hension (chatgpt) du probl\u00e8me des tables dans la bdd\"\n$targetIdx = 0\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $idx++\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $anchorText) {\n        $targetIdx = $idx\n    }\n}\nif ($targetIdx -eq 0) {\n    throw \"Anchor paragraph not found\"\n}\n\n$target = $d.Paragraphs.Item($targetIdx)\n\n# 1) \"25/07 :\" -- top level (ilvl 0 / ListLevelNumber 1)\n$target.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Item($targetIdx + 1)\n$p1.Range.InsertBefore(\"25/07 :\")\n$p1.Range.ListFormat.ListLevelNumber = 1\n\n# 2) \"Faits :\" -- second level (ilvl 1 / ListLevelNumber 2)\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($targetIdx + 2)\n$p2.Range.InsertBefore(\"Faits :\")\n$p2.Range.ListFormat.ListLevelNumber = 2\n\n# 3) The new fact itself -- third level (ilvl 2 / ListLevelNumber 3)\n$p2.Range.InsertParagraphAfter()\n$p3 = $d.Paragraphs.Item($targetIdx + 3)\n$p3.Range.InsertBefore(\"R\u00e9solution de la double table dans la bdd mais re-probl\u00e8me avec les mdp (.matches())\")\n$p3.Range.ListFormat.ListLevelNumber = 3\n"}
